# The deck originally has 10 slides:
#   1  Fall 2021 Precalc / Lesson 6.3
#   2  Do now...
#   3  (framing)
#   4  B24 rules
#   5  Independent work...
#   6  Today's activity: practice problems
#   7  Problem (1d)              <- removed
#   8  Problem (2b)               <- removed
#   9  Reflection                 <- removed
#   10 wrapping up! (be sure to...) <- becomes the new slide 7
#
# Slides 7, 8 and 9 ("Problem (1d)", "Problem (2b)" and "Reflection") are
# deleted, leaving the former slide 10 ("wrapping up!") as the new, final
# slide 7 of a 7-slide deck.

$p = $ppt.ActivePresentation

# Deleting the slide currently at index 7 three times removes the three
# unwanted slides in turn (each deletion shifts the following slides up by
# one position), leaving the old slide 10 as the new slide 7.
$p.Slides.Item(7).Delete()
$p.Slides.Item(7).Delete()
$p.Slides.Item(7).Delete()
